# "union casos de uso"
#
# Row 7 used to hold the (now removed) "ADMIN / Crear alumno / No ser pudo
# guardar alumno..." use case. All the rows from 7 downward get rewritten:
# the old rows 8-10 (Docente use cases) shift up to 7-9, a new "Eliminar
# tema / predecesor" row becomes row 10, and four brand-new rows (11-14)
# are appended before the existing blank filler rows. The now-superfluous
# trailing blank row 52 is removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: Docente / Carga/Edita nuevo concepto / agregar actividades a un concepto deshabilitado
$ws.Cells.Item(7, 1).Value = "Docente"
$ws.Cells.Item(7, 2).Value = "Carga/Edita nuevo concepto"
$ws.Cells.Item(7, 3).Value = "agregar actividades a un concepto deshabilitado"
$ws.Rows.Item(7).RowHeight = 30

# --- Row 8: Docente / eliminar concepto / Se auto relacione las depencias cuando se elimina un concepto o tema
$ws.Cells.Item(8, 1).Value = "Docente"
$ws.Cells.Item(8, 2).Value = "eliminar concepto"
$ws.Cells.Item(8, 3).Value = "Se auto relacione las depencias cuando se elimina un concepto o tema"
$ws.Rows.Item(8).AutoFit()

# --- Row 9: Docente / editar actividad / get de la actividad pincha cuando le damos editar
$ws.Cells.Item(9, 1).Value = "Docente"
$ws.Cells.Item(9, 2).Value = "editar actividad"
$ws.Cells.Item(9, 3).Value = "get de la actividad pincha cuando le damos editar"

# --- Row 10: Docente / Eliminar tema / cuando el tema está eliminado no deberia ...
$ws.Cells.Item(10, 1).Value = "Docente"
$ws.Cells.Item(10, 2).Value = "Eliminar tema"
$ws.Cells.Item(10, 3).Value = "cuando el tema está eliminado no deberia aparecer para poder elegirse como predecesor."

# --- Row 11: docente / Seleccionar predecesora / Selección que quede encolumnada
$ws.Cells.Item(11, 1).Value = "docente"
$ws.Cells.Item(11, 2).Value = "Seleccionar predecesora"
$ws.Cells.Item(11, 3).Value = "Selección que quede encolumnada"
$ws.Rows.Item(11).RowHeight = 30

# --- Row 12: alumno / renovar contraseña / no se renueva la contraseña desde el admin
$ws.Cells.Item(12, 1).Value = "alumno"
$ws.Cells.Item(12, 2).Value = "renovar contraseña"
$ws.Cells.Item(12, 3).Value = "no se renueva la contraseña desde el admin "
$ws.Rows.Item(12).RowHeight = 30

# --- Row 13: alumno / boton volver / cambiar el boton volver de temas y conceptos ...
$ws.Cells.Item(13, 1).Value = "alumno"
$ws.Cells.Item(13, 2).Value = "boton volver"
$ws.Cells.Item(13, 3).Value = "cambiar el boton volver de temas y conceptos cambiar la posición del mismo arriba a la derecha"
$ws.Rows.Item(13).RowHeight = 30

# --- Row 14: docente / carga y modificacion de temas y conceptos / Estan invertidos ...
$ws.Cells.Item(14, 1).Value = "docente"
$ws.Cells.Item(14, 2).Value = "carga y modificacion de temas y conceptos"
$ws.Cells.Item(14, 3).Value = "Estan invertidos los campos copete y descripción"
$ws.Rows.Item(14).RowHeight = 45

# Drop the now-unused trailing blank row (was row 52, dimension shrinks to D51)
$ws.Rows.Item(52).Delete()

# Update the view: scrolled to row 2, active cell/selection on C15
$ws.Activate()
$ws.Range("C15").Select()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
